$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.518.11'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.62%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.811.05'
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.008'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.71%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.006'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '308.00'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.82%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4551'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.95%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3661'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.38%  '
$ws.Range("E9").Value = '  -2.22%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8767'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.21%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07778'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.02%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '19.36'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.68%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.805.00'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.45%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.282'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.09%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.359'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.31%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '86.26'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.29%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.008'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.67%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008590'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.74%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.005'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.41%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.580.16'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.48%  '
$ws.Range("E21").Value = '  -2.98%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.57%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.45'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.987'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.95%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.54'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.24%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '17.92'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.55%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.059'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.73%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '112.93'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.64%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.860'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.41%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08679'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.065'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.45%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.528'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.11%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7335'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.118'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.670'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.85%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.005'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.84%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.083'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.36%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01947'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.42%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05114'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.98%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.907'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.986'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.70%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5008'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.27%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1560'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.12%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.162'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.83%  '
$ws.Range("E45").Value = '  +0.65%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4605'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.96%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.970'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.22%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '100.89'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.85%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.588'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.35%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06005'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.21%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '64.48'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.53%  '
